$d = $word.ActiveDocument

$d.Content.Find.Execute("22×95=2090", $true, $false, $false, $false, $false, $true, 1, $false, "56×56=3136", 2)
$d.Content.Find.Execute("71×74=5254", $true, $false, $false, $false, $false, $true, 1, $false, "60×11=660", 2)
$d.Content.Find.Execute("25×37=925", $true, $false, $false, $false, $false, $true, 1, $false, "71×31=2201", 2)
$d.Content.Find.Execute("60×79=4740", $true, $false, $false, $false, $false, $true, 1, $false, "22×75=1650", 2)
$d.Content.Find.Execute("30×36=1080", $true, $false, $false, $false, $false, $true, 1, $false, "43×18=774", 2)
$d.Content.Find.Execute("72×74=5328", $true, $false, $false, $false, $false, $true, 1, $false, "64×54=3456", 2)
$d.Content.Find.Execute("29×85=2465", $true, $false, $false, $false, $false, $true, 1, $false, "78×73=5694", 2)
$d.Content.Find.Execute("16×36=576", $true, $false, $false, $false, $false, $true, 1, $false, "34×26=884", 2)
$d.Content.Find.Execute("56×81=4536", $true, $false, $false, $false, $false, $true, 1, $false, "26×22=572", 2)
$d.Content.Find.Execute("80×29=2320", $true, $false, $false, $false, $false, $true, 1, $false, "33×29=957", 2)
$d.Content.Find.Execute("29×96=2784", $true, $false, $false, $false, $false, $true, 1, $false, "26×31=806", 2)
$d.Content.Find.Execute("50×86=4300", $true, $false, $false, $false, $false, $true, 1, $false, "85×65=5525", 2)
$d.Content.Find.Execute("83×25=2075", $true, $false, $false, $false, $false, $true, 1, $false, "40×23=920", 2)
$d.Content.Find.Execute("52×62=3224", $true, $false, $false, $false, $false, $true, 1, $false, "88×24=2112", 2)
$d.Content.Find.Execute("98×85=8330", $true, $false, $false, $false, $false, $true, 1, $false, "14×66=924", 2)
$d.Content.Find.Execute("32×86=2752", $true, $false, $false, $false, $false, $true, 1, $false, "92×67=6164", 2)
$d.Content.Find.Execute("47×68=3196", $true, $false, $false, $false, $false, $true, 1, $false, "15×96=1440", 2)
$d.Content.Find.Execute("43×33=1419", $true, $false, $false, $false, $false, $true, 1, $false, "99×19=1881", 2)
$d.Content.Find.Execute("81×70=5670", $true, $false, $false, $false, $false, $true, 1, $false, "61×15=915", 2)
$d.Content.Find.Execute("76×31=2356", $true, $false, $false, $false, $false, $true, 1, $false, "21×66=1386", 2)
$d.Content.Find.Execute("41×91=3731", $true, $false, $false, $false, $false, $true, 1, $false, "79×36=2844", 2)
$d.Content.Find.Execute("50×78=3900", $true, $false, $false, $false, $false, $true, 1, $false, "92×50=4600", 2)
$d.Content.Find.Execute("49×74=3626", $true, $false, $false, $false, $false, $true, 1, $false, "72×26=1872", 2)
$d.Content.Find.Execute("24×42=1008", $true, $false, $false, $false, $false, $true, 1, $false, "93×65=6045", 2)
$d.Content.Find.Execute("72×22=1584", $true, $false, $false, $false, $false, $true, 1, $false, "64×57=3648", 2)
